# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-158) from 2023-09-21 (serial 45190) to 2023-09-23 (serial 45192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C158").Value2 = 45192
